# Fix Training Data Issue (#48)
# Data was taken from 1 day off due to way NBA stats were shown.
#
# Column BF ("Date") held the literal text "5-14-2007-08" for every data
# row (rows 2-31). Correct it to the literal text "2008-05-14".
#
# Simply assigning the date-like text via .Value/.Value2/.Formula makes
# Excel's smart cell-entry parser coerce the string into a real date
# serial number (and attach a date number format/style to the cell) -
# that is standard Excel behaviour, but not what we want here: the
# original (and desired) cell is a plain text cell with no special
# style. To land literal text without Excel's autoconvert kicking in,
# we compute the string via a formula, then copy/paste-special just the
# computed value back over itself (xlPasteValues), which keeps the cell
# as plain text and leaves formatting untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteValues = -4163
$newDate = "2008-05-14"

for ($r = 2; $r -le 31; $r++) {
    $cell = $ws.Range("BF$r")
    $cell.Formula = '="' + $newDate + '"'
    $cell.Copy()
    $cell.PasteSpecial($xlPasteValues)
}

$excel.CutCopyMode = 0
